$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.363.25"
$ws.Range("D3").Value = "2.912.18"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.03"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.91"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "3.373.38"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("D17").Value = "2.910.15"
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").Value = "52.387.58"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.98"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.36"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.173"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("E30").Value = "  +8.92%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.37"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +12.94%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0984"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.18"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +6.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +13.63%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +7.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.37%  "
$ws.Range("D48").Value = "2.199.84"
$ws.Range("E48").Value = "  +3.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.265"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +22.96%  "
$ws.Range("E50").Value = "  +11.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.969"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.94%  "
